# Weekly update: insert a new price record for "Haba" at
# Terminal Hortofrutícola Agro Chillán.
#
# The new observation is inserted as row 100 of the data table, pushing
# the previously-existing rows 100-104 down to become rows 101-105
# (dimension grows from A1:R104 to A1:R105).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 100, shifting rows 100:104 down to 101:105.
$ws.Rows.Item(100).Insert()

# Populate the newly inserted row 100 with this week's data.
$ws.Cells.Item(100, 1).Value  = 7
$ws.Cells.Item(100, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(100, 3).Value  = "Ñuble"
$ws.Cells.Item(100, 4).Value  = 45267
$ws.Cells.Item(100, 5).Value  = 16
$ws.Cells.Item(100, 6).Value  = 100112026
$ws.Cells.Item(100, 7).Value  = "Haba"
$ws.Cells.Item(100, 8).Value  = "Sin especificar"
$ws.Cells.Item(100, 9).Value  = "Primera"
$ws.Cells.Item(100, 10).Value = 100
$ws.Cells.Item(100, 11).Value = 10000
$ws.Cells.Item(100, 12).Value = 12000
$ws.Cells.Item(100, 13).Value = 11000
$ws.Cells.Item(100, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(100, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(100, 16).Value = 440
$ws.Cells.Item(100, 17).Value = 25
$ws.Cells.Item(100, 18).Value = "Hortaliza"
